# NMRA_Region_Division_Map.xlsx — add a new division row.
#
# A new division ("Allegheny Western Division", division #12 under region
# #32 / "MCR") is inserted as row 148 of the "Division Reassignments"
# sheet; every row that used to be 148-201 shifts down by one (149-202).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Division Reassignments")

# Push row 148 (and everything below it) down by one row, carrying the
# existing number/text styles (s="1"/s="2") onto the freshly inserted row.
$ws.Rows.Item(148).Insert()

# Populate the newly-opened row 148 with the new division's data.
$ws.Range("A148").Value = 32
$ws.Range("B148").Value = 12
$ws.Range("C148").Value = "Allegheny Western Division"
$ws.Range("D148").Value = "MCR"

# Match the author's final cursor position/selection.
$ws.Range("B148").Select() | Out-Null
